$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry: same task/ID as row 2, re-run finished at 16:03 taking 10 min.
$taskId = $ws.Range("A2").Value2
$status = $ws.Range("D2").Value2

$ws.Range("A3").Value = $taskId
$ws.Range("B3").Value = "19.12.2025 16:03"
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = $status
